$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I20").Value = 0.349748557851123
$ws.Range("J20").Value = 0.2242144557103363
$ws.Range("K20").Value = 0.1948034178860719
$ws.Range("L20").Value = 2.540403574172122
